# Add columns I (I0) and J (IF) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - uses the same style as the other header cells (s="1")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data rows 2-19 for columns I and J
$values = @(
    @(7, 7),
    @(12, 13),
    @(6, 7),
    @(8, 8),
    @(11, 11),
    @(5, 7),
    @(9, 9),
    @(8, 8),
    @(4, 6),
    @(4, 6),
    @(3, 5),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(6, 6),
    @(6, 7),
    @(8, 8),
    @(8, 8)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}

$wb.Save()
